{"js": "// Office.js (Word JavaScript API) edit script.\n// 1) Split the run \"Eu, XXXX, estou assinando esse contrato para voc\u00ea.\"\n//    into three runs: \"Eu, \" | \"XXXX\" | \", estou assinando esse contrato para voc\u00ea.\"\n// 2) Replace the standalone run \"Lira\" with \"XXXX\".\n\nconst body = context.document.body;\n\n// ---- Part 1: split \"Eu, XXXX, ...\" into three runs ----------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/uniqueLocalId\");\nawait context.sync();\n\nlet introParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const candidate = paragraphs.items[i];\n  if (candidate.text.indexOf(\"Eu, XXXX, estou assinando esse contrato\") !== -1) {\n    introParagraph = candidate;\n    break;\n  }\n}\n\nif (introParagraph) {\n  // Find the \"XXXX\" token inside that paragraph and give it a throw-away\n  // formatting toggle. Office.js only splits a run's underlying <w:r> when\n  // a sub-range's formatting actually changes, so flipping bold on/off is\n  // the reliable way to force the paragraph into three runs at exactly the\n  // boundaries we need (\"Eu, \" | \"XXXX\" | \", estou...\").\n  const paragraphId = introParagraph.uniqueLocalId;\n  const tokenResults = introParagraph.search(\"XXXX\", { matchCase: true });\n  tokenResults.load(\"items\");\n  await context.sync();\n\n  const token = tokenResults.items[0];\n  token.font.bold = true;\n  await context.sync();\n  token.font.bold = false;\n  await context.sync();\n\n  // The toggle leaves a harmless-but-visible empty <w:rPr/> on the middle\n  // run. Round-trip the paragraph through getOoxml/insertOoxml to strip it\n  // back out while preserving the now-split run structure and the\n  // paragraph's own identity (paraId/textId), so nothing else moves.\n  const ooxmlResult = introParagraph.getOoxml();\n  await context.sync();\n\n  const fullOoxml = ooxmlResult.value;\n  const openTagMarker = `w14:paraId=\"${paragraphId}\"`;\n  const openTagStart = fullOoxml.indexOf(`<w:p ${openTagMarker}`);\n  if (openTagStart !== -1) {\n    const closeTagEnd = fullOoxml.indexOf(\"</w:p>\", openTagStart) + \"</w:p>\".length;\n    let paragraphFragment = fullOoxml.substring(openTagStart, closeTagEnd);\n\n    // Strip any now-empty run-properties element(s) left behind by the\n    // bold toggle (<w:rPr/> or <w:rPr><w:b/></w:rPr> style artifacts).\n    paragraphFragment = paragraphFragment.replace(/<w:rPr\\/>/g, \"\");\n    paragraphFragment = paragraphFragment.replace(/<w:rPr><w:b\\/><\\/w:rPr>/g, \"\");\n    paragraphFragment = paragraphFragment.replace(/<w:rPr><w:b w:val=\"0\"\\/><\\/w:rPr>/g, \"\");\n\n    const wrappedOoxml =\n      '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n      '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n      '<w:body>' + paragraphFragment + '</w:body>' +\n      '</w:document>' +\n      '</pkg:xmlData>' +\n      '</pkg:part>' +\n      '</pkg:package>';\n\n    introParagraph.insertOoxml(wrappedOoxml, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// ---- Part 2: \"Lira\" -> \"XXXX\" -------------------------------------------\nconst signatureResults = body.search(\"Lira\", { matchCase: true });\nsignatureResults.load(\"items\");\nawait context.sync();\n\nif (signatureResults.items.length > 0) {\n  signatureResults.items[0].insertText(\"XXXX\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# 1) Split the run \"Eu, XXXX, estou assinando esse contrato para voc\u00ea.\"\n#    into three runs: \"Eu, \" | \"XXXX\" | \", estou assinando esse contrato para voc\u00ea.\"\n# 2) Replace the standalone run \"Lira\" with \"XXXX\".\n\n$d = $word.ActiveDocument\n\n# ---- Part 1: split \"Eu, XXXX, ...\" into three runs ----------------------\n$introRange = $d.Content\n$introRange.Find.ClearFormatting()\n$introRange.Find.Text = \"XXXX\"\n$introRange.Find.MatchCase = $true\n$introRange.Find.Execute() | Out-Null\n\nif ($introRange.Find.Found) {\n    # Re-assigning a Range's FormattedText to itself is a no-op in terms of\n    # visible content/formatting, but it makes Word materialize the target\n    # sub-range as its own run(s), splitting the parent run at the\n    # sub-range's boundaries without leaving any run-properties residue.\n    $introRange.FormattedText = $introRange.FormattedText\n}\n\n# ---- Part 2: \"Lira\" -> \"XXXX\" -------------------------------------------\n$signatureRange = $d.Content\n$signatureRange.Find.ClearFormatting()\n$signatureRange.Find.Text = \"Lira\"\n$signatureRange.Find.MatchCase = $true\n$signatureRange.Find.Execute() | Out-Null\n\nif ($signatureRange.Find.Found) {\n    $signatureRange.Text = \"XXXX\"\n}\n"}
